$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 82, pushing existing rows 82-106 down to 84-108.
$ws.Rows("82:83").Insert()

# Row 82: new data row (Jengibre, Primera)
$ws.Range("A82").Value = 9
$ws.Range("B82").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C82").Value = "Metropolitana"
$ws.Range("D82").Value = 44809
$ws.Range("E82").Value = 13
$ws.Range("F82").Value = 100114007
$ws.Range("G82").Value = "Jengibre"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 350
$ws.Range("K82").Value = 13000
$ws.Range("L82").Value = 15000
$ws.Range("M82").Value = 14143
$ws.Range("N82").Value = "$/caja 13 kilos"
$ws.Range("O82").Value = "Perú"
$ws.Range("P82").Value = 1088
$ws.Range("Q82").Value = 13
$ws.Range("R82").Value = "Hortaliza"

# Row 83: new data row (Jengibre, Segunda)
$ws.Range("A83").Value = 9
$ws.Range("B83").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C83").Value = "Metropolitana"
$ws.Range("D83").Value = 44809
$ws.Range("E83").Value = 13
$ws.Range("F83").Value = 100114007
$ws.Range("G83").Value = "Jengibre"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Segunda"
$ws.Range("J83").Value = 160
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = 10000
$ws.Range("N83").Value = "$/caja 13 kilos"
$ws.Range("O83").Value = "Perú"
$ws.Range("P83").Value = 769
$ws.Range("Q83").Value = 13
$ws.Range("R83").Value = "Hortaliza"
